$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049779759218478
$ws.Range("D2").Value = 1.046719066332265
$ws.Range("E2").Value = 1.055985252378274
$ws.Range("F2").Value = 1.064647919388419
$ws.Range("I2").Value = 1.037891377006614
$ws.Range("J2").Value = 1.054816396963882
$ws.Range("K2").Value = 1.049483715523722
$ws.Range("L2").Value = 1.058724206674621
$ws.Range("M2").Value = 1.067363290219444

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051422459474212
$ws.Range("D3").Value = 1.047919464890361
$ws.Range("E3").Value = 1.057491640126352
$ws.Range("F3").Value = 1.066341427089236
$ws.Range("I3").Value = 1.038281166959797
$ws.Range("J3").Value = 1.056105125808102
$ws.Range("K3").Value = 1.05049480413219
$ws.Range("L3").Value = 1.060042379528147
$ws.Range("M3").Value = 1.06886984643475

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052483056015862
$ws.Range("D4").Value = 1.048693893065423
$ws.Range("E4").Value = 1.058464461607183
$ws.Range("F4").Value = 1.06743553210545
$ws.Range("I4").Value = 1.038530648265959
$ws.Range("J4").Value = 1.056936347906306
$ws.Range("K4").Value = 1.051146174450459
$ws.Range("L4").Value = 1.060892912643394
$ws.Range("M4").Value = 1.069842515340154

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052928382075461
$ws.Range("D5").Value = 1.049018917043821
$ws.Range("E5").Value = 1.058872987815108
$ws.Range("F5").Value = 1.067895096139295
$ws.Range("I5").Value = 1.038634877674889
$ws.Range("J5").Value = 1.057285163016612
$ws.Range("K5").Value = 1.051419330001377
$ws.Range("L5").Value = 1.061249908183964
$ws.Range("M5").Value = 1.070250915126877

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053003122437076
$ws.Range("D6").Value = 1.049073458217405
$ws.Range("E6").Value = 1.058941555179441
$ws.Range("F6").Value = 1.067972235961908
$ws.Range("I6").Value = 1.038652340077506
$ws.Range("J6").Value = 1.05734369388453
$ws.Range("K6").Value = 1.051465154314117
$ws.Range("L6").Value = 1.061309816226209
$ws.Range("M6").Value = 1.070319457622475

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052489008624616
$ws.Range("D7").Value = 1.048698238185375
$ws.Range("E7").Value = 1.058469922100866
$ws.Range("F7").Value = 1.067441674369937
$ws.Range("I7").Value = 1.038532043543998
$ws.Range("J7").Value = 1.056941011256555
$ws.Range("K7").Value = 1.051149827032658
$ws.Range("L7").Value = 1.060897685056856
$ws.Range("M7").Value = 1.069847974388837

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050335409035239
$ws.Range("D8").Value = 1.047125229401512
$ws.Range("E8").Value = 1.056494745022865
$ws.Range("F8").Value = 1.065220608550963
$ws.Range("I8").Value = 1.038023677354945
$ws.Range("J8").Value = 1.055252487491488
$ws.Range("K8").Value = 1.049826016437621
$ws.Range("L8").Value = 1.059170194438135
$ws.Range("M8").Value = 1.067872894169934

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046522046270962
$ws.Range("D9").Value = 1.04433536220758
$ws.Range("E9").Value = 1.052999143638482
$ws.Range("F9").Value = 1.061293227389028
$ws.Range("I9").Value = 1.037106751968268
$ws.Range("J9").Value = 1.052256228533951
$ws.Range("K9").Value = 1.047470986582578
$ws.Range("L9").Value = 1.056107247268729
$ws.Range("M9").Value = 1.06437544449139

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043966671542212
$ws.Range("D10").Value = 1.042462862256721
$ws.Range("E10").Value = 1.0506579772415
$ws.Range("F10").Value = 1.05866514036275
$ws.Range("I10").Value = 1.036481066908012
$ws.Range("J10").Value = 1.050244127970701
$ws.Range("K10").Value = 1.045885527599279
$ws.Range("L10").Value = 1.054052009273657
$ws.Range("M10").Value = 1.062031666704188

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042856887425303
$ws.Range("D11").Value = 1.041648957671885
$ws.Range("E11").Value = 1.04964153364423
$ws.Range("F11").Value = 1.057524660118352
$ws.Range("I11").Value = 1.036206676841026
$ws.Range("J11").Value = 1.049369277515474
$ws.Range("K11").Value = 1.045195245158145
$ws.Range("L11").Value = 1.053158791693067
$ws.Range("M11").Value = 1.06101376037572

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042444155707767
$ws.Range("D12").Value = 1.041346162354365
$ws.Range("E12").Value = 1.049263563034532
$ws.Range("F12").Value = 1.057100645555483
$ws.Range("I12").Value = 1.036104231654054
$ws.Range("J12").Value = 1.049043767548334
$ws.Range("K12").Value = 1.04493826872574
$ws.Range("L12").Value = 1.052826506058557
$ws.Range("M12").Value = 1.060635196146535

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042532711301675
$ws.Range("D13").Value = 1.041411134617791
$ws.Range("E13").Value = 1.049344658133247
$ws.Range("F13").Value = 1.057191615961116
$ws.Range("I13").Value = 1.036126230294261
$ws.Range("J13").Value = 1.049113615678971
$ws.Range("K13").Value = 1.044993417193805
$ws.Range("L13").Value = 1.052897805481501
$ws.Range("M13").Value = 1.060716420861421

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042822781348343
$ws.Range("D14").Value = 1.041623938258286
$ws.Range("E14").Value = 1.049610299077271
$ws.Range("F14").Value = 1.057489618965766
$ws.Range("I14").Value = 1.036198219413001
$ws.Range("J14").Value = 1.04934238207596
$ws.Range("K14").Value = 1.04517401518227
$ws.Range("L14").Value = 1.0531313352228
$ws.Range("M14").Value = 1.060982477768454

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043001435271376
$ws.Range("D15").Value = 1.041754990396265
$ws.Range("E15").Value = 1.049773913450695
$ws.Range("F15").Value = 1.057673176439892
$ws.Range("I15").Value = 1.036242504662504
$ws.Range("J15").Value = 1.049483259184322
$ws.Range("K15").Value = 1.045285211136156
$ws.Range("L15").Value = 1.053275153369448
$ws.Range("M15").Value = 1.061146341753089

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044040253733996
$ws.Range("D16").Value = 1.042516812331289
$ws.Range("E16").Value = 1.050725377275117
$ws.Range("F16").Value = 1.058740776421532
$ws.Range("I16").Value = 1.036499203941615
$ws.Range("J16").Value = 1.050302112175511
$ws.Range("K16").Value = 1.045931259241768
$ws.Range("L16").Value = 1.05411121899995
$ws.Range("M16").Value = 1.062099156784269

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044690986306283
$ws.Range("D17").Value = 1.042993846741078
$ws.Range("E17").Value = 1.0513214732015
$ws.Range("F17").Value = 1.059409774450105
$ws.Range("I17").Value = 1.036659294297538
$ws.Range("J17").Value = 1.050814786068749
$ws.Range("K17").Value = 1.046335493199459
$ws.Range("L17").Value = 1.054634773281222
$ws.Range("M17").Value = 1.062696011111115

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045070231135889
$ws.Range("D18").Value = 1.043271794388822
$ws.Range("E18").Value = 1.051668905883736
$ws.Range("F18").Value = 1.05979974927885
$ws.Range("I18").Value = 1.036752338313101
$ws.Range("J18").Value = 1.051113473710273
$ws.Range("K18").Value = 1.046570913038084
$ws.Range("L18").Value = 1.054939837348562
$ws.Range("M18").Value = 1.06304385415257

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045199490570629
$ws.Range("D19").Value = 1.043366517008089
$ws.Range("E19").Value = 1.051787327715471
$ws.Range("F19").Value = 1.059932680251582
$ws.Range("I19").Value = 1.03678400738995
$ws.Range("J19").Value = 1.051215260059889
$ws.Range("K19").Value = 1.046651123810352
$ws.Range("L19").Value = 1.055043803000408
$ws.Range("M19").Value = 1.063162410490578

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044621201696359
$ws.Range("D20").Value = 1.042942696400482
$ws.Range("E20").Value = 1.051257544729013
$ws.Range("F20").Value = 1.059338022257551
$ws.Range("I20").Value = 1.036642152691698
$ws.Range("J20").Value = 1.050759816918103
$ws.Range("K20").Value = 1.046292160329786
$ws.Range("L20").Value = 1.054578633657542
$ws.Range("M20").Value = 1.062632004582675

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042737377100783
$ws.Range("D21").Value = 1.041561286074654
$ws.Range("E21").Value = 1.049532086082276
$ws.Range("F21").Value = 1.057401875368559
$ws.Range("I21").Value = 1.036177034923129
$ws.Range("J21").Value = 1.049275031388863
$ws.Range("K21").Value = 1.045120849513276
$ws.Range("L21").Value = 1.053062580561339
$ws.Range("M21").Value = 1.060904143674525

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041549992742401
$ws.Range("D22").Value = 1.04068998622039
$ws.Range("E22").Value = 1.048444796090099
$ws.Range("F22").Value = 1.05618228305237
$ws.Range("I22").Value = 1.035881560302274
$ws.Range("J22").Value = 1.048338291107296
$ws.Range("K22").Value = 1.044381069942319
$ws.Range("L22").Value = 1.05210645125612
$ws.Range("M22").Value = 1.059815052077416

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042179731971779
$ws.Range("D23").Value = 1.041152142752546
$ws.Range("E23").Value = 1.049021423079611
$ws.Range("F23").Value = 1.056829030818246
$ws.Range("I23").Value = 1.036038486182639
$ws.Range("J23").Value = 1.048835181583012
$ws.Range("K23").Value = 1.044773559650081
$ws.Range("L23").Value = 1.052613594655007
$ws.Range("M23").Value = 1.060392661876031

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044652735335578
$ws.Range("D24").Value = 1.042965809959946
$ws.Range("E24").Value = 1.051286432057331
$ws.Range("F24").Value = 1.059370444727667
$ws.Range("I24").Value = 1.036649899277541
$ws.Range("J24").Value = 1.050784656181427
$ws.Range("K24").Value = 1.046311741710012
$ws.Range("L24").Value = 1.0546040017143
$ws.Range("M24").Value = 1.062660927279985

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04751015266267
$ws.Range("D25").Value = 1.045058793808688
$ws.Range("E25").Value = 1.053904691922899
$ws.Range("F25").Value = 1.062310230743898
$ws.Range("I25").Value = 1.037346322542319
$ws.Range("J25").Value = 1.053033363803286
$ws.Range("K25").Value = 1.048082506830726
$ws.Range("L25").Value = 1.056901389603565
$ws.Range("M25").Value = 1.06528171208973
